$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-10 (TruckID, AssignedDockPosition, start_loading_time, end_loading_time)
$data = @(
    @(3, 1, 5, 5),
    @(7, 1, 10, 10),
    @(8, 1, 15, 15),
    @(2, 2, 5, 5),
    @(4, 2, 10, 10),
    @(5, 2, 15, 17),
    @(8, 2, 22, 22),
    @(1, 3, 5, 5),
    @(6, 3, 10, 10)
)

$row = 2
foreach ($rowData in $data) {
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $row++
}
